$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182567834854126
$ws.Range("B1").Value = 2.330402135848999
$ws.Range("C1").Value = 3.763520240783691
$ws.Range("D1").Value = 3.128907203674316
$ws.Range("E1").Value = 1.1434725522995
